# Updates the cryptocurrency price/volume snapshot table (columns D and E,
# rows 2-51) to the latest scraped values.
#
# Note: several "Price" values (column D) are plain decimal numbers as text
# (e.g. "603.74"). Assigning them straight to .Value would make Excel
# auto-detect them as numbers (changing cell type/formatting). To keep them
# as plain text with the original (default) cell style, we briefly switch
# the cell to a Text number format, assign the value, then restore the
# original "Normal" style so no stray formatting is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.937.25'
$ws.Range("E2").Value = '  -1.19%  '
$ws.Range("D3").Value = '3.142.35'
$ws.Range("E3").Value = '  -0.60%  '
$ws.Range("E4").Value = '  -0.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '603.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.03%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.26%  '
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.139.31'
$ws.Range("E8").Value = '  -0.62%  '
$ws.Range("E9").Value = '  -0.55%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.150'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.78%  '
$ws.Range("E11").Value = '  -2.44%  '
$ws.Range("E12").Value = '  -1.53%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000254'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  -2.00%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.99'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -2.56%  '
$ws.Range("D15").Value = '3.659.86'
$ws.Range("E15").Value = '  -0.52%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.121'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +2.81%  '
$ws.Range("D17").Value = '63.977.60'
$ws.Range("D18").Value = '3.144.67'
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("E19").Value = '  -0.94%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '487.91'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.02%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.71'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.711'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -1.23%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.71'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -3.16%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '88.29'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.51%  '
$ws.Range("E25").Value = '  -3.61%  '
$ws.Range("E26").Value = '  -0.02%  '
$ws.Range("E27").Value = '  -2.40%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '8.21'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -4.45%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.00'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.55%  '
$ws.Range("E30").Value = '  -1.06%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.62'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.112'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.18%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.66'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.69%  '
$ws.Range("E35").Value = '  -2.61%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.06'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.59%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '52.73'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.60%  '
$ws.Range("D38").Value = '0.0₃0747'
$ws.Range("E38").Value = '  -4.20%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.98'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -7.34%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '434.88'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -6.10%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0398'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.66%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.121'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.08%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '8.34'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.83%  '
$ws.Range("D44").Value = '2.936.05'
$ws.Range("E44").Value = '  +3.04%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.261'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -2.78%  '
$ws.Range("E46").Value = '  -5.15%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.41'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.81%  '
$ws.Range("E48").Value = '  -0.11%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '25.87'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -2.84%  '
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '120.33'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.05%  '
